$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add a "Mesa" (table) column header in C1 ---
# C1 already carries style s="2" (inherited from the header row/col
# defaults) so a plain value assignment keeps that style untouched.
$ws.Range("C1").Value = "Mesa"

# --- Row 2: guest renamed, phone renumbered, table number added ---
$ws.Range("A2").Value = "Mouse"

# B2 must stay a *text* value ("+549456"), not get auto-coerced to a
# number by losing the leading "+". Force text via NumberFormat, then
# strip the formatting change back out (copy the plain format from A2,
# a cell with no explicit style) so no stray "s" attribute is left
# behind on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "+549456"
$ws.Cells.Item(2, 1).Copy()
$ws.Range("B2").PasteSpecial(-4122)

# New C2 cell (table number, numeric). Newly created cells in these
# columns would otherwise inherit the sheet's column-level default
# style (s="2"); copy formats from A2 (unstyled) to keep it plain.
$ws.Range("C2").Value = 1
$ws.Cells.Item(2, 1).Copy()
$ws.Range("C2").PasteSpecial(-4122)

# --- Row 3: guest renamed, phone renumbered, "N/A" table marker added ---
$ws.Range("A3").Value = "Paolo"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "+549"
$ws.Cells.Item(3, 1).Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "N/A"
$ws.Cells.Item(3, 1).Copy()
$ws.Range("C3").PasteSpecial(-4122)

# --- Remove old row 4 (duplicate "Paolo Cetti" entry) entirely ---
$ws.Rows("4").Delete()

$excel.CutCopyMode = $false
